$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-500: bump the date value from 45203 to 45204
$range = $ws.Range("C2:C500")
$range.Value = 45204
